# Generate Report for Handback
#
# Row 6 in both the "zh-cn" and "de-de" sheets corresponds to the
# 166616ad-d08f-49c8-b815-9f1115236e60 source file, whose handback was
# picked up by the report generator. The generator now:
#   - fills in the "Latest Target File" (I), linking it back to the
#     handoff file on GitHub,
#   - records the produced "Latest Handback File" (J),
#   - stamps the "Latest Handback DateTime" (K),
#   - flags that the handback isn't based on the latest handoff via the
#     "Error Detail" column (P).
# The report's Error Detail column is also widened so the message is
# readable.

$wb = $excel.ActiveWorkbook

$handoffTarget = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/466d7c936b2ea6b52aa6b98d3592275c3ee3d318/e2e/166616ad-d08f-49c8-b815-9f1115236e60.md"
$handoffDisplay = "166616ad-d08f-49c8-b815-9f1115236e60.md"

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e648dc8b20b99471df10fe01aea29498dfabbc51/e2e/166616ad-d08f-49c8-b815-9f1115236e60.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/466d7c936b2ea6b52aa6b98d3592275c3ee3d318/e2e/166616ad-d08f-49c8-b815-9f1115236e60.md."

# ---- zh-cn sheet ----
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("J6").Value = "166616ad-d08f-49c8-b815-9f1115236e60.977cffc582c0f6ada02e0057eb5d330a6f96c928.zh-cn.xlf"
$wsZhCn.Range("K6").Value = "2016-08-31 14:54:46"
$wsZhCn.Range("P6").Value = $errorDetail

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I6"), $handoffTarget, [System.Type]::Missing, [System.Type]::Missing, $handoffDisplay)

$wsZhCn.Columns.Item(16).ColumnWidth = 40

# ---- de-de sheet ----
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("J6").Value = "166616ad-d08f-49c8-b815-9f1115236e60.977cffc582c0f6ada02e0057eb5d330a6f96c928.de-de.xlf"
$wsDeDe.Range("K6").Value = "2016-08-31 14:54:54"
$wsDeDe.Range("P6").Value = $errorDetail

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I6"), $handoffTarget, [System.Type]::Missing, [System.Type]::Missing, $handoffDisplay)

$wsDeDe.Columns.Item(16).ColumnWidth = 40
